$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.Value = "'" + $value
    $cell.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") '43.342.57'
Set-TextValue $ws.Range("E2") '  -2.08%  '

# Row 3
Set-TextValue $ws.Range("D3") '2.239.05'
Set-TextValue $ws.Range("E3") '  -2.05%  '

# Row 4
Set-TextValue $ws.Range("E4") '  -0.08%  '

# Row 5
Set-TextValue $ws.Range("D5") '230.57'
Set-TextValue $ws.Range("E5") '  -1.01%  '

# Row 6
Set-TextValue $ws.Range("D6") '0.639'
Set-TextValue $ws.Range("E6") '  -0.22%  '

# Row 7
Set-TextValue $ws.Range("D7") '63.78'
Set-TextValue $ws.Range("E7") '  -2.11%  '

# Row 8
Set-TextValue $ws.Range("E8") '  +0.05%  '

# Row 9
Set-TextValue $ws.Range("D9") '0.441'
Set-TextValue $ws.Range("E9") '  +1.38%  '

# Row 10
Set-TextValue $ws.Range("D10") '0.0953'
Set-TextValue $ws.Range("E10") '  -8.68%  '

# Row 11
Set-TextValue $ws.Range("D11") '56.46'
Set-TextValue $ws.Range("E11") '  -0.01%  '

# Row 12
Set-TextValue $ws.Range("D12") '27.62'
Set-TextValue $ws.Range("E12") '  +5.96%  '

# Row 13
Set-TextValue $ws.Range("E13") '  -1.52%  '

# Row 14
Set-TextValue $ws.Range("D14") '2.571.28'
Set-TextValue $ws.Range("E14") '  -2.07%  '

# Row 15
Set-TextValue $ws.Range("D15") '15.42'
Set-TextValue $ws.Range("E15") '  -4.15%  '

# Row 16
Set-TextValue $ws.Range("D16") '6.06'
Set-TextValue $ws.Range("E16") '  +0.27%  '

# Row 17
Set-TextValue $ws.Range("D17") '0.825'
Set-TextValue $ws.Range("E17") '  -1.54%  '

# Row 18
Set-TextValue $ws.Range("D18") '2.242.46'
Set-TextValue $ws.Range("E18") '  -1.96%  '

# Row 19
Set-TextValue $ws.Range("D19") '43.220.98'
Set-TextValue $ws.Range("E19") '  -2.04%  '

# Row 20
Set-TextValue $ws.Range("D20") '0.0₃0964'
Set-TextValue $ws.Range("E20") '  -4.17%  '

# Row 21
Set-TextValue $ws.Range("D21") '72.79'
Set-TextValue $ws.Range("E21") '  -2.01%  '

# Row 22
Set-TextValue $ws.Range("D22") '6.08'
Set-TextValue $ws.Range("E22") '  -0.67%  '

# Row 23
Set-TextValue $ws.Range("D23") '245.77'
Set-TextValue $ws.Range("E23") '  -6.17%  '

# Row 24
Set-TextValue $ws.Range("D24") '1.00'
Set-TextValue $ws.Range("E24") '  -0.02%  '

# Row 25
Set-TextValue $ws.Range("D25") '3.66'
Set-TextValue $ws.Range("E25") '  +29.67%  '

# Row 26
Set-TextValue $ws.Range("D26") '2.40'
Set-TextValue $ws.Range("E26") '  -3.52%  '

# Row 27
Set-TextValue $ws.Range("D27") '2.21'
Set-TextValue $ws.Range("E27") '  -4.77%  '

# Row 28
Set-TextValue $ws.Range("D28") '9.73'
Set-TextValue $ws.Range("E28") '  -4.44%  '

# Row 29
Set-TextValue $ws.Range("D29") '173.10'
Set-TextValue $ws.Range("E29") '  +0.79%  '

# Row 30
Set-TextValue $ws.Range("D30") '21.48'
Set-TextValue $ws.Range("E30") '  +1.77%  '

# Row 31
$ws.Range("B31").Value = 'Kaspa'
$ws.Range("C31").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue $ws.Range("D31") '0.128'
Set-TextValue $ws.Range("E31") '  -7.40%  '

# Row 32
$ws.Range("B32").Value = 'ImmutableX'
$ws.Range("C32").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue $ws.Range("D32") '1.40'
Set-TextValue $ws.Range("E32") '  -2.26%  '

# Row 33
Set-TextValue $ws.Range("E33") '  -0.37%  '

# Row 34
Set-TextValue $ws.Range("D34") '4.92'
Set-TextValue $ws.Range("E34") '  +3.35%  '

# Row 35
Set-TextValue $ws.Range("D35") '0.0673'
Set-TextValue $ws.Range("E35") '  -2.34%  '

# Row 36
Set-TextValue $ws.Range("D36") '4.88'
Set-TextValue $ws.Range("E36") '  -2.57%  '

# Row 37
Set-TextValue $ws.Range("D37") '3.59'
Set-TextValue $ws.Range("E37") '  -7.40%  '

# Row 38
Set-TextValue $ws.Range("D38") '6.27'
Set-TextValue $ws.Range("E38") '  -8.47%  '

# Row 39
Set-TextValue $ws.Range("E39") '  -4.11%  '

# Row 40
Set-TextValue $ws.Range("D40") '0.0250'
Set-TextValue $ws.Range("E40") '  +0.02%  '

# Row 41
Set-TextValue $ws.Range("D41") '0.999'
Set-TextValue $ws.Range("E41") '  -0.18%  '

# Row 42
Set-TextValue $ws.Range("D42") '8.64'
Set-TextValue $ws.Range("E42") '  +1.20%  '

# Row 43
Set-TextValue $ws.Range("D43") '4.46'
Set-TextValue $ws.Range("E43") '  -0.09%  '

# Row 44
Set-TextValue $ws.Range("D44") '16.96'
Set-TextValue $ws.Range("E44") '  -4.51%  '

# Row 45
Set-TextValue $ws.Range("D45") '96.29'
Set-TextValue $ws.Range("E45") '  -2.23%  '

# Row 46
$ws.Range("B46").Value = 'Cronos'
$ws.Range("C46").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue $ws.Range("D46") '0.0940'
Set-TextValue $ws.Range("E46") '  -3.33%  '

# Row 47
$ws.Range("B47").Value = 'TerraClassic'
$ws.Range("C47").Value = 'https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc'
Set-TextValue $ws.Range("D47") '0.000210'
Set-TextValue $ws.Range("E47") '  -0.03%  '

# Row 48
Set-TextValue $ws.Range("E48") '  -2.80%  '

# Row 49
Set-TextValue $ws.Range("D49") '1.439.59'
Set-TextValue $ws.Range("E49") '  -2.48%  '

# Row 50
$ws.Range("B50").Value = 'Celestia'
$ws.Range("C50").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
Set-TextValue $ws.Range("D50") '9.84'
Set-TextValue $ws.Range("E50") '  +0.49%  '

# Row 51
$ws.Range("B51").Value = 'NEARProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue $ws.Range("D51") '2.27'
Set-TextValue $ws.Range("E51") '  -3.42%  '

